{"js": "// The site footer block that used to follow the \"Requisitos\" section\n// (\"Ver no Jupiter Salvar em pdf Salvar em docx\" and the \"\u00a9 2020 ...\"\n// copyright line, plus the blank paragraph that separated them from the\n// requirement text) was removed when the site was rebuilt. Delete those\n// paragraphs, leaving the trailing blank paragraph (and page break) intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\nlet startIndex = -1;\nlet endIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (text === targets[0]) {\n    startIndex = i;\n  } else if (text === targets[1] && startIndex !== -1) {\n    endIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1 && endIndex !== -1) {\n  // Also remove the blank paragraph immediately preceding the \"Ver no\n  // Jupiter...\" paragraph, which separated it from the requirements text.\n  let from = startIndex;\n  if (from - 1 >= 0 && items[from - 1].text.trim() === \"\") {\n    from = from - 1;\n  }\n  for (let i = endIndex; i >= from; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# The site footer block that used to follow the \"Requisitos\" section\n# (\"Ver no Jupiter Salvar em pdf Salvar em docx\" and the \"\u00a9 2020 ...\"\n# copyright line, plus the blank paragraph that separated them from the\n# requirement text) was removed when the site was rebuilt. Delete those\n# paragraphs, leaving the trailing blank paragraph (and page break) intact.\n\n$d = $word.ActiveDocument\n\n$jupiterMarker = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightMarker = \"Powered by Jekyll and Github pages\"\n\n$startPara = $null\n$endPara = $null\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text\n    if ($text -like \"*$jupiterMarker*\") {\n        $startPara = $i\n    } elseif (($text -like \"*$copyrightMarker*\") -and ($startPara -ne $null)) {\n        $endPara = $i\n        break\n    }\n}\n\nif (($startPara -ne $null) -and ($endPara -ne $null)) {\n    # Also remove the blank paragraph immediately preceding the \"Ver no\n    # Jupiter...\" paragraph, which separated it from the requirements text.\n    $fromPara = $startPara\n    if ($fromPara -gt 1) {\n        $prevText = $d.Paragraphs.Item($fromPara - 1).Range.Text.Trim()\n        if ($prevText -eq \"\") {\n            $fromPara = $fromPara - 1\n        }\n    }\n\n    $rangeStart = $d.Paragraphs.Item($fromPara).Range.Start\n    $rangeEnd = $d.Paragraphs.Item($endPara).Range.End\n\n    $deleteRange = $d.Range($rangeStart, $rangeEnd)\n    $deleteRange.Delete()\n}\n"}
